$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell D1 ("success"): clone an existing header cell's
# formatting (bold/centered/bordered look shared by B1 "list" and C1
# "count") onto D1, then overwrite its text.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "success"

# Fill the success indicator for rows 2-10: "1" for the row whose count is
# the large/successful one (383, row 6), "0" for every other row.
#
# Writing the numeric-looking strings "0"/"1" straight into .Value would be
# auto-coerced to numbers by Excel. Instead, enter them as formulas that
# evaluate to text ("0"/"1"), then convert the whole range to values via
# copy / paste-special-values - this yields genuine text cells (stored as
# shared strings) without leaving behind any formulas or extra cell
# formatting/number-format styles.
for ($r = 2; $r -le 10; $r++) {
    $count = $ws.Cells.Item($r, 3).Value()
    if ($count -eq 383) {
        $ws.Cells.Item($r, 4).Formula = "=""1"""
    } else {
        $ws.Cells.Item($r, 4).Formula = "=""0"""
    }
}

$successRange = $ws.Range("D2:D10")
$successRange.Copy()
$successRange.PasteSpecial(-4163)
